# Refresh the cryptos price/volume snapshot (GitHub Actions data pull).
# Column D ("Price") values are stored as text that merely looks numeric
# (e.g. "68.184.87"), so a leading apostrophe is used to keep Excel from
# auto-converting them to real numbers when the cell's .Value is set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''68.184.87'
$ws.Range('E2').Value = '  +1.37%  '
$ws.Range('D3').Value = '''3.276.51'
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''586.84'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').Value = '''184.93'
$ws.Range('E6').Value = '  +3.54%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '''0.601'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +4.26%  '
$ws.Range('D10').Value = '''6.72'
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('D12').Value = '''3.844.65'
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('E13').Value = '  +0.53%  '
$ws.Range('D14').Value = '''28.64'
$ws.Range('E14').Value = '  +2.49%  '
$ws.Range('D15').Value = '''68.194.69'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('E16').Value = '  +2.11%  '
$ws.Range('D17').Value = '''3.272.92'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').Value = '''13.64'
$ws.Range('E19').Value = '  +2.52%  '
$ws.Range('D20').Value = '''382.54'
$ws.Range('E20').Value = '  +2.48%  '
$ws.Range('E21').Value = '  +2.40%  '
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('E24').Value = '  +1.27%  '
$ws.Range('E25').Value = '  +1.95%  '
$ws.Range('D26').Value = '''0.193'
$ws.Range('E26').Value = '  +7.29%  '
$ws.Range('D27').Value = '''9.78'
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').Value = '''1.00'
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('D29').Value = '''5.82'
$ws.Range('E29').Value = '  +3.33%  '
$ws.Range('E30').Value = '  +0.93%  '
$ws.Range('D31').Value = '''23.00'
$ws.Range('E31').Value = '  +1.95%  '
$ws.Range('D32').Value = '''7.22'
$ws.Range('E32').Value = '  +5.97%  '
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').Value = '  +2.83%  '
$ws.Range('E36').Value = '  +0.61%  '
$ws.Range('E37').Value = '  +0.97%  '
$ws.Range('D38').Value = '''0.839'
$ws.Range('E38').Value = '  -2.29%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''6.75'
$ws.Range('E39').Value = '  -0.66%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').Value = '''26.65'
$ws.Range('E40').Value = '  -0.59%  '
$ws.Range('E41').Value = '  +4.96%  '
$ws.Range('E42').Value = '  +0.66%  '
$ws.Range('D43').Value = '''25.63'
$ws.Range('E43').Value = '  -0.46%  '
$ws.Range('E44').Value = '  +2.46%  '
$ws.Range('D45').Value = '''41.27'
$ws.Range('E45').Value = '  +2.06%  '
$ws.Range('D46').Value = '''2.630.68'
$ws.Range('E46').Value = '  -4.94%  '
$ws.Range('D47').Value = '''342.57'
$ws.Range('E47').Value = '  -4.29%  '
$ws.Range('E48').Value = '  +2.37%  '
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('E51').Value = '  -0.09%  '
